$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K" = strikeouts),
# replacing the previous "Strike#" counts with regenerated strikeout totals.
$gValues = @{
    2 = 2
    3 = 0
    4 = 0
    5 = 0
    6 = 2
    7 = 1
    8 = 1
    9 = 1
    10 = 0
    11 = 0
    12 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 3
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 2
    45 = 2
    46 = 1
    47 = 2
    49 = 3
    50 = 1
    51 = 3
    52 = 0
    53 = 1
    54 = 0
    55 = 0
    56 = 0
    57 = 3
    59 = 1
    61 = 1
    62 = 1
    63 = 2
    64 = 2
    65 = 1
    66 = 1
    67 = 1
    68 = 1
    69 = 2
    70 = 0
    71 = 1
    72 = 1
    73 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
